$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.295.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.643.52'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.58%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.61%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.643.19'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.61%  '
$ws.Range('E10').Value = '  +8.52%  '
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('E13').Value = '  +2.03%  '
$ws.Range('E14').Value = '  +3.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.96'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.123.63'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.504.63'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.645.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '364.96'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.46'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('E22').Value = '  +3.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.87'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.55%  '
$ws.Range('E24').Value = '  -0.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '75.49'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -1.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000106'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.30%  '
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '561.55'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.87%  '
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('E34').Value = '  +0.98%  '
$ws.Range('E35').Value = '  +3.06%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  +3.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '161.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.33%  '
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.375'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.00%  '
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.36'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0₆0339'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.64'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('E45').Value = '  +2.08%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '156.45'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.63%  '
$ws.Range('E49').Value = '  +1.84%  '
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.80'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.68%  '
